# Aadhaar_Extracted_Data_A9.xlsx -- correct the extracted record in row 2.
# The previously extracted name ("Anoop Negi") was wrong; it belongs to
# "surishta Devi". The address/state split was also off: "garhwal" actually
# belongs to the address line, and the State column should just say
# "Uttarakhand".
#
# NOTE: order matters here -- new shared-string entries are appended in the
# order cells are written, and must match the order the strings were
# authored in upstream (State, then Address, then Name) to reproduce the
# canonical shared-strings table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = "Uttarakhand"
$ws.Range("E2").Value = "w/o Devendra Singh  00.00 jhandi chaur  pachhimi pauri  garhwal "
$ws.Range("A2").Value = "surishta Devi"

# Column widths, tuned for the new (longer) name/DOB columns and the
# now-wider address column (matches Excel's auto-fit result for these
# columns after the text changed above).
$ws.Columns.Item(1).ColumnWidth = 18.333333333333336
$ws.Columns.Item(2).ColumnWidth = 20.833333333333336
$ws.Columns.Item(5).ColumnWidth = 50.666666666666664

# Leave the selection where it was left after the edit.
$ws.Range("F7").Select()
